# PlanningEffectif.xlsx - "Grosse modification de la doc" edit
#
# Fills in the missing hours for the "Analyse" (row 9), "Implémentation"
# (row 10) and "Documentation" (row 13) lines for the second week
# (columns H:K = "Semaine du 14 au 20 février") and updates the active
# selection, matching the author's manual edit + format-painter pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Analyse (row 9): Mercredi (I9) now logs 1h00 -----------------------
# Paint the same "worked hours" look already used on H9/G9 (time format +
# light-blue fill) before typing the value in, exactly like copying the
# format from a neighbouring filled-in cell.
$ws.Range("H9").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$ws.Range("I9").Value = 0.041666666666666664

# --- Implementation (row 10): Mercredi (I10) 7h00, Jeudi (J10) 2h30 -----
$ws.Range("H10").Copy()
$ws.Range("I10").PasteSpecial(-4122)
$ws.Range("I10").Value = 0.2916666666666667

$ws.Range("H10").Copy()
$ws.Range("J10").PasteSpecial(-4122)
$ws.Range("J10").Value = 0.10416666666666667

# --- Documentation (row 13): Jeudi (J13) 5h30 ----------------------------
# Copy the format already applied on G13 (bold font + light-blue fill) so
# the new entry matches the rest of that row's "logged" cells.
$ws.Range("G13").Copy()
$ws.Range("J13").PasteSpecial(-4122)
$ws.Range("J13").Value = 0.22916666666666666

$excel.CutCopyMode = 0

# --- Leftover selection state left by the author on save ---------------
$ws.Range("I18").Select()
